$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")
$ws.Range("A2").Value = "juan.perez99_65904@test.com"
